# Inserts a new week of "Apio" price data (date 45041) at rows 987:988 of
# the single worksheet, pushing the existing rows 987-1036 down to 989-1038.
# This mirrors the upstream diff: dimension grows from A1:R1036 to A1:R1038
# and two brand-new rows (Primera / Segunda quality) are added for the new
# reporting date, while every previously existing row keeps its original
# values (just shifted down by two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 987:1036 down by two rows, leaving 987:988 blank.
$ws.Rows("987:988").Insert()

# --- New row 987 ("Primera" quality) -----------------------------------
$ws.Cells.Item(987, 1).Value = 6
$ws.Cells.Item(987, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(987, 3).Value = "Metropolitana"
$ws.Cells.Item(987, 4).Value2 = 45041
$ws.Cells.Item(987, 5).Value = 13
$ws.Cells.Item(987, 6).Value = 100112017
$ws.Cells.Item(987, 7).Value = "Apio"
$ws.Cells.Item(987, 8).Value = "Americana (o)"
$ws.Cells.Item(987, 9).Value = "Primera"
$ws.Cells.Item(987, 10).Value = 2800
$ws.Cells.Item(987, 11).Value = 7000
$ws.Cells.Item(987, 12).Value = 8000
$ws.Cells.Item(987, 13).Value = 7571
$ws.Cells.Item(987, 14).Value = "$/docena de matas"
$ws.Cells.Item(987, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(987, 16).Value = 1262
$ws.Cells.Item(987, 17).Value = 6
$ws.Cells.Item(987, 18).Value = "Hortaliza"

# --- New row 988 ("Segunda" quality) ------------------------------------
$ws.Cells.Item(988, 1).Value = 6
$ws.Cells.Item(988, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(988, 3).Value = "Metropolitana"
$ws.Cells.Item(988, 4).Value2 = 45041
$ws.Cells.Item(988, 5).Value = 13
$ws.Cells.Item(988, 6).Value = 100112017
$ws.Cells.Item(988, 7).Value = "Apio"
$ws.Cells.Item(988, 8).Value = "Americana (o)"
$ws.Cells.Item(988, 9).Value = "Segunda"
$ws.Cells.Item(988, 10).Value = 1200
$ws.Cells.Item(988, 11).Value = 5000
$ws.Cells.Item(988, 12).Value = 5000
$ws.Cells.Item(988, 13).Value = 5000
$ws.Cells.Item(988, 14).Value = "$/docena de matas"
$ws.Cells.Item(988, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(988, 16).Value = 833
$ws.Cells.Item(988, 17).Value = 6
$ws.Cells.Item(988, 18).Value = "Hortaliza"
